$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.966.90"
$ws.Range("E2").Value = "  -0.70%  "

$ws.Range("D3").Value = "3.065.50"
$ws.Range("E3").Value = "  -3.29%  "

$ws.Range("D5").Value = "'587.78"
$ws.Range("E5").Value = "  +0.15%  "

$ws.Range("D6").Value = "'130.85"
$ws.Range("E6").Value = "  -3.10%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "3.062.37"
$ws.Range("E8").Value = "  -3.43%  "

$ws.Range("E9").Value = "  -0.08%  "

$ws.Range("E10").Value = "  -1.09%  "

$ws.Range("D11").Value = "'5.25"
$ws.Range("E11").Value = "  +0.11%  "

$ws.Range("E12").Value = "  -2.03%  "

$ws.Range("D13").Value = "'0.0000237"
$ws.Range("E13").Value = "  +1.21%  "

$ws.Range("D14").Value = "'33.83"
$ws.Range("E14").Value = "  +2.12%  "

$ws.Range("E15").Value = "  +0.71%  "

$ws.Range("D16").Value = "3.572.61"
$ws.Range("E16").Value = "  -2.85%  "

$ws.Range("D17").Value = "62.243.27"
$ws.Range("E17").Value = "  -0.25%  "

$ws.Range("D18").Value = "3.069.46"
$ws.Range("E18").Value = "  -3.29%  "

$ws.Range("D19").Value = "'6.41"
$ws.Range("E19").Value = "  -1.66%  "

$ws.Range("D20").Value = "'449.32"
$ws.Range("E20").Value = "  -1.20%  "

$ws.Range("E21").Value = "  -1.84%  "

$ws.Range("E22").Value = "  -3.91%  "

$ws.Range("E23").Value = "  -2.78%  "

$ws.Range("D24").Value = "'12.94"
$ws.Range("E24").Value = "  -3.18%  "

$ws.Range("D25").Value = "'80.99"
$ws.Range("E25").Value = "  -2.60%  "

$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("E27").Value = "  -0.06%  "

$ws.Range("E28").Value = "  -3.41%  "

$ws.Range("D29").Value = "'2.02"
$ws.Range("E29").Value = "  -0.03%  "

$ws.Range("D30").Value = "'7.47"
$ws.Range("E30").Value = "  -3.73%  "

$ws.Range("D31").Value = "'6.46"
$ws.Range("E31").Value = "  -6.11%  "

$ws.Range("D32").Value = "'26.03"
$ws.Range("E32").Value = "  -4.66%  "

$ws.Range("D33").Value = "'0.0978"
$ws.Range("E33").Value = "  -6.00%  "

$ws.Range("D34").Value = "'2.33"
$ws.Range("E34").Value = "  -2.50%  "

$ws.Range("D35").Value = "'0.976"
$ws.Range("E35").Value = "  -5.35%  "

$ws.Range("D36").Value = "'5.75"
$ws.Range("E36").Value = "  -2.02%  "

$ws.Range("D37").Value = "'50.57"
$ws.Range("E37").Value = "  -1.13%  "

$ws.Range("D38").Value = "0.0₃0697"
$ws.Range("E38").Value = "  +0.93%  "

$ws.Range("D39").Value = "'0.0377"
$ws.Range("E39").Value = "  -1.96%  "

$ws.Range("E40").Value = "  -0.18%  "

$ws.Range("D41").Value = "'0.109"
$ws.Range("E41").Value = "  -1.97%  "

$ws.Range("D42").Value = "'382.43"
$ws.Range("E42").Value = "  -6.93%  "

$ws.Range("D43").Value = "'2.53"
$ws.Range("E43").Value = "  -4.85%  "

$ws.Range("D44").Value = "2.705.50"
$ws.Range("E44").Value = "  -5.96%  "

$ws.Range("D46").Value = "'124.65"
$ws.Range("E46").Value = "  +0.04%  "

$ws.Range("D47").Value = "'0.241"
$ws.Range("E47").Value = "  -3.00%  "

$ws.Range("D48").Value = "'2.03"
$ws.Range("E48").Value = "  -4.71%  "

$ws.Range("D49").Value = "'34.18"
$ws.Range("E49").Value = "  -5.48%  "

$ws.Range("E50").Value = "  -1.44%  "

$ws.Range("D51").Value = "'24.14"
$ws.Range("E51").Value = "  -4.45%  "
